$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The B column currently holds text like "28000 [21000 - 36000]".
# Replace it with just the leading numeric estimate as a real number.
$values = @(28000, 8500, 0, 580, 23000, 46000, 13000, 6100, 150000, 910, 38000, 6100, 3600, 19000, 7800, 72000, 53000, 240000, 48000, 38000)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

$ws.Range("B22").Select() | Out-Null

$wb.Save()
